$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 209.2
$ws.Range("J38").Value = 429.5
$ws.Range("L38").Value = 1288.5
$ws.Range("N38").Value = -2032.5
$ws.Range("H58").Value = 1618.8889
$ws.Range("J58").Value = 2035.7142
$ws.Range("L58").Value = 6107.142599999999
$ws.Range("N58").Value = -6407.142599999999
$ws.Range("H64").Value = 5051.5293
$ws.Range("I64").Value = 3624.2856
$ws.Range("J64").Value = 6050.6
$ws.Range("K64").Value = 3624.2856
$ws.Range("L64").Value = 6050.6
$ws.Range("M64").Value = -3376.2856
$ws.Range("N64").Value = -6546.6
$ws.Range("H67").Value = 5051.5293
$ws.Range("I67").Value = 3624.2856
$ws.Range("J67").Value = 6050.6
$ws.Range("K67").Value = 3624.2856
$ws.Range("L67").Value = 6050.6
$ws.Range("M67").Value = -2766.2856
$ws.Range("N67").Value = -7766.6
$ws.Range("H87").Value = 80000
$ws.Range("J87").Value = 80000
$ws.Range("L87").Value = 80000
$ws.Range("N87").Value = -82496
$ws.Range("H90").Value = 80000
$ws.Range("J90").Value = 80000
$ws.Range("L90").Value = 240000
$ws.Range("N90").Value = -252480
$ws.Range("H137").Value = 1027.4348
$ws.Range("I137").Value = 904.2222
$ws.Range("K137").Value = 2712.6666
$ws.Range("M137").Value = -162.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 26712.95
$ws.Range("I2").Value = 35498.62
$ws.Range("J2").Value = 1234.5
$ws.Range("K2").Value = 35498.62
$ws.Range("L2").Value = 1234.5
$ws.Range("M2").Value = -35385.62
$ws.Range("N2").Value = -1460.5
$ws.Range("H61").Value = 2605.111
$ws.Range("I61").Value = 2371.1177
$ws.Range("K61").Value = 2371.1177
$ws.Range("M61").Value = -2159.1177
$ws.Range("H74").Value = 6676.0415
$ws.Range("I74").Value = 1349.6923
$ws.Range("K74").Value = 1349.6923
$ws.Range("M74").Value = -475.6922999999999
$ws.Range("H77").Value = 6676.0415
$ws.Range("I77").Value = 1349.6923
$ws.Range("K77").Value = 6748.461499999999
$ws.Range("M77").Value = -2380.461499999999
$ws.Range("H116").Value = 26712.95
$ws.Range("I116").Value = 35498.62
$ws.Range("J116").Value = 1234.5
$ws.Range("K116").Value = 35498.62
$ws.Range("L116").Value = 1234.5
$ws.Range("M116").Value = -33204.62
$ws.Range("N116").Value = -5822.5
$ws.Range("H136").Value = 2605.111
$ws.Range("I136").Value = 2371.1177
$ws.Range("K136").Value = 7113.353099999999
$ws.Range("M136").Value = -4563.353099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 26712.95
$ws.Range("I3").Value = 35498.62
$ws.Range("J3").Value = 1234.5
$ws.Range("K3").Value = 35498.62
$ws.Range("L3").Value = 1234.5
$ws.Range("M3").Value = -35384.62
$ws.Range("N3").Value = -1462.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 45457708
$ws.Range("I86").Value = 100002296
$ws.Range("J86").Value = 3879.5
$ws.Range("K86").Value = 100002296
$ws.Range("L86").Value = 3879.5
$ws.Range("M86").Value = -100001173
$ws.Range("N86").Value = -6125.5
$ws.Range("H89").Value = 45457708
$ws.Range("I89").Value = 100002296
$ws.Range("J89").Value = 3879.5
$ws.Range("K89").Value = 500011480
$ws.Range("L89").Value = 19397.5
$ws.Range("M89").Value = -500005864
$ws.Range("N89").Value = -30629.5
$ws.Range("H99").Value = 10430685
$ws.Range("I99").Value = 12516582
$ws.Range("J99").Value = 1200
$ws.Range("K99").Value = 12516582
$ws.Range("L99").Value = 1200
$ws.Range("M99").Value = -12515084
$ws.Range("N99").Value = -4196
$ws.Range("H122").Value = 1437.5
$ws.Range("I122").Value = 1280
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 3840
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -1390
$ws.Range("N122").Value = -10000
$ws.Range("H126").Value = 10430685
$ws.Range("I126").Value = 12516582
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 37549746
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -37547276
$ws.Range("N126").Value = -8540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1139.2
$ws.Range("I8").Value = 1139.2
$ws.Range("K8").Value = 3417.6
$ws.Range("M8").Value = -3278.6
$ws.Range("H122").Value = 1099.8
$ws.Range("I122").Value = 248.4
$ws.Range("J122").Value = 1951.2
$ws.Range("K122").Value = 2235.6
$ws.Range("L122").Value = 17560.8
$ws.Range("M122").Value = 214.4000000000001
$ws.Range("N122").Value = -22460.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1491.4762
$ws.Range("I16").Value = 1354.7333
$ws.Range("J16").Value = 1833.3334
$ws.Range("K16").Value = 1354.7333
$ws.Range("L16").Value = 1833.3334
$ws.Range("M16").Value = -1184.7333
$ws.Range("N16").Value = -2173.3334
$ws.Range("H46").Value = 789
$ws.Range("I46").Value = 722.5
$ws.Range("J46").Value = 833.3333
$ws.Range("K46").Value = 722.5
$ws.Range("L46").Value = 833.3333
$ws.Range("M46").Value = -534.5
$ws.Range("N46").Value = -1209.3333
$ws.Range("H55").Value = 410.5
$ws.Range("I55").Value = 380.75
$ws.Range("J55").Value = 430.33334
$ws.Range("K55").Value = 380.75
$ws.Range("L55").Value = 430.33334
$ws.Range("M55").Value = -207.75
$ws.Range("N55").Value = -776.33334
$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1251
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6256
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 3294.9565
$ws.Range("I122").Value = 2048.5
$ws.Range("K122").Value = 6145.5
$ws.Range("M122").Value = -3695.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5298
$ws.Range("H62").Value = 16647.666
$ws.Range("I62").Value = 20200.2
$ws.Range("J62").Value = 13418.091
$ws.Range("K62").Value = 20200.2
$ws.Range("L62").Value = 13418.091
$ws.Range("M62").Value = -19576.2
$ws.Range("N62").Value = -14666.091
$ws.Range("H65").Value = 16647.666
$ws.Range("I65").Value = 20200.2
$ws.Range("J65").Value = 13418.091
$ws.Range("K65").Value = 101001
$ws.Range("L65").Value = 67090.455
$ws.Range("M65").Value = -97881
$ws.Range("N65").Value = -73330.455
$ws.Range("H81").Value = 3576.842
$ws.Range("I81").Value = 2084.4443
$ws.Range("J81").Value = 4920
$ws.Range("K81").Value = 4168.8886
$ws.Range("L81").Value = 9840
$ws.Range("M81").Value = -3107.8886
$ws.Range("N81").Value = -11962
$ws.Range("H84").Value = 3576.842
$ws.Range("I84").Value = 2084.4443
$ws.Range("J84").Value = 4920
$ws.Range("K84").Value = 20844.443
$ws.Range("L84").Value = 49200
$ws.Range("M84").Value = -15540.443
$ws.Range("N84").Value = -59808
$ws.Range("H122").Value = 112433.445
$ws.Range("I122").Value = 126237.625
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 378712.875
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -376262.875
$ws.Range("N122").Value = -10900
$ws.Range("H126").Value = 250925
$ws.Range("I126").Value = 250925
$ws.Range("K126").Value = 752775
$ws.Range("M126").Value = -750305
